# Update the "Historia de Usuario" paragraph:
#  - Shorten "Como ... usuario del sistema de gestión contravencional" down to
#    "Como ... usuario del sistema " and turn the manual line break that used
#    to follow it into a real paragraph break (so "Quiero ..." starts its own
#    paragraph instead of just being a new line in the same paragraph).
#  - Remove the now-orphaned manual line break run that used to sit right
#    before "Quiero".
#  - Give both the "Como ..." paragraph and the new "Quiero ... Para ..."
#    paragraph spacing-after = 0 (w:spacing w:after="0").

$d = $word.ActiveDocument

# Step 1: trim the run text and turn the trailing manual line break into a
# paragraph mark, splitting the paragraph in two right after "sistema ".
$rFind = $d.Content
$found = $rFind.Find.Execute(
    "usuario del sistema de gestión contravencional", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "usuario del sistema ^p", 2)
if (-not $found) {
    throw "Could not find the 'usuario del sistema de gestión contravencional' text to replace"
}

# Step 2: the paragraph split leaves the original manual line break (the one
# that used to separate "...contravencional" from "Quiero") dangling right
# before "Quiero" in the new paragraph. Find it (^l matches a manual line
# break) and delete just that single character, preserving the "Quiero" run
# (and its bold formatting) untouched.
$rBreak = $d.Content
$foundBreak = $rBreak.Find.Execute(
    "^lQuiero", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundBreak) {
    throw "Could not find the leftover manual line break before 'Quiero'"
}
$rBreak.End = $rBreak.Start + 1
$rBreak.Delete()

# Step 3: set spacing-after = 0 on both the "Como ..." paragraph and the new
# "Quiero ... Para ..." paragraph that resulted from the split.
$pComo = $d.Paragraphs.Item(2)
$pComo.Format.SpaceAfter = 0

$pQuiero = $d.Paragraphs.Item(3)
$pQuiero.Format.SpaceAfter = 0
